$wb = $excel.ActiveWorkbook

# Sheet 2: Ciserano Italy
$ws = $wb.Worksheets.Item(2)
$ws.Range("H4").Value = 0.012
$ws.Range("I4").Value = 0.0122
$ws.Range("J4").Value = 0.0242
$ws.Range("K4").Value = 0.0122
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("E5").Value = 0.866666666666667
$ws.Range("E6").Value = 0.866666666666667
$ws.Range("E7").Value = 0.866666666666667
$ws.Range("L7").Value = 0.9231
$ws.Range("M7").Value = 0.866666666666667
$ws.Range("N7").Value = 0.866666666666667
$ws.Range("O7").Value = 0.866666666666667
$ws.Range("P7").Value = 0.866666666666667
$ws.Range("Q7").Value = 0.866666666666667
$ws.Range("R7").Value = 0.866666666666667
$ws.Range("S7").Value = 0.866666666666667
$ws.Range("T7").Value = 0.866666666666667
$ws.Range("U7").Value = 0.866666666666667
$ws.Range("V7").Value = 0.866666666666667
$ws.Range("W7").Value = 0.866666666666667
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776

# Sheet 7: Gainsborough - Epc
$ws = $wb.Worksheets.Item(7)
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# Sheet 10: Molndal Sweden
$ws = $wb.Worksheets.Item(10)
$ws.Range("G5").Value = $null
$ws.Range("J5").Value = $null

# Sheet 11: Nove Mesto Slovakia
$ws = $wb.Worksheets.Item(11)
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# Sheet 12: Piedras Negras Fasco Mexico
$ws = $wb.Worksheets.Item(12)
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("L7").Value = $null

# Sheet 13: Rotherham United Kingdom
$ws = $wb.Worksheets.Item(13)
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# Sheet 14: Sao Paulo Brazil
$ws = $wb.Worksheets.Item(14)
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# Sheet 16: Waldenburg Germany
$ws = $wb.Worksheets.Item(16)
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("G9").Value = 0.0291
$ws.Range("J9").Value = 0.0293
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0

# Sheet 18: Juarez MEJ II
$ws = $wb.Worksheets.Item(18)
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# Sheet 19: Yueyang China
$ws = $wb.Worksheets.Item(19)
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# Sheet 21: Changzhou Epc China
$ws = $wb.Worksheets.Item(21)
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("L6").Value = $null

# Sheet 23: Juarez Casa II
$ws = $wb.Worksheets.Item(23)
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# Sheet 24: Bangalore India
$ws = $wb.Worksheets.Item(24)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("L4").Value = 0.5

# Sheet 25: Black River Falls Wisconsin
$ws = $wb.Worksheets.Item(25)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776

# Sheet 26: Juarez FCDM
$ws = $wb.Worksheets.Item(26)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("L4").Value = $null

# Sheet 27: Mumbai India
$ws = $wb.Worksheets.Item(27)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 0.0833333333333333
$ws.Range("N4").Value = 0.25
$ws.Range("O4").Value = 0.0833333333333333
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# Sheet 28: Noida India
$ws = $wb.Worksheets.Item(28)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("G4").Value = 0.0222
$ws.Range("I4").Value = 0.0227
$ws.Range("J4").Value = 0.0447
$ws.Range("K4").Value = 0.0227
$ws.Range("L4").Value = 0.0488

# Sheet 29: Juarez Casa I
$ws = $wb.Worksheets.Item(29)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("L4").Value = $null

# Sheet 30: Suzhou China
$ws = $wb.Worksheets.Item(30)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# Sheet 31: Edmonton EDM Canada
$ws = $wb.Worksheets.Item(31)
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
